# Scheduled-runner update: refresh currentAveragePrice* / LevePrice* /
# LeveProfit* figures (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# leve-profit tables, row by row, to the latest market-board snapshot.
$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "ALC"; Row = 9;   Cells = @{ "H" = 237.25;     "I" = 83;         "J" = 700;        "K" = 83;         "L" = 700;        "M" = 86;               "N" = -1038 } },
    @{ Sheet = "ALC"; Row = 95;  Cells = @{ "H" = 35624;      "J" = 35624;      "L" = 35624;       "N" = -41116 } },
    @{ Sheet = "ALC"; Row = 97;  Cells = @{ "H" = 1004.6842;  "I" = 0;          "K" = 0;           "M" = $null } },
    @{ Sheet = "ALC"; Row = 129; Cells = @{ "H" = 1110.9487;  "J" = 1332.2333;  "L" = 3996.699900000001; "N" = -13996.6999 } },

    @{ Sheet = "ARM"; Row = 41;  Cells = @{ "H" = 3685.3333;  "I" = 3685.3333;  "J" = 0;           "K" = 3685.3333;  "L" = 0;          "M" = -3271.3333;       "N" = $null } },
    @{ Sheet = "ARM"; Row = 45;  Cells = @{ "H" = 2589.862;   "I" = 2195.3333;  "K" = 2195.3333;    "M" = -1818.3333 } },
    @{ Sheet = "ARM"; Row = 61;  Cells = @{ "H" = 2099.6897;  "I" = 1575.64;    "K" = 1575.64;      "M" = -1363.64 } },
    @{ Sheet = "ARM"; Row = 74;  Cells = @{ "H" = 66669096;   "I" = 100001130; "J" = 5022.8;       "K" = 100001130; "L" = 5022.8;      "M" = -100000256;       "N" = -6770.8 } },
    @{ Sheet = "ARM"; Row = 77;  Cells = @{ "H" = 66669096;   "I" = 100001130; "J" = 5022.8;       "K" = 500005650; "L" = 25114;       "M" = -500001282;       "N" = -33850 } },
    @{ Sheet = "ARM"; Row = 132; Cells = @{ "H" = 16461.795;  "I" = 1377.0834;  "J" = 52665.1;     "K" = 4131.2502;  "L" = 157995.3;   "M" = -1601.2502;       "N" = -163055.3 } },
    @{ Sheet = "ARM"; Row = 136; Cells = @{ "H" = 2099.6897;  "I" = 1575.64;    "K" = 4726.92;      "M" = -2176.92 } },

    @{ Sheet = "BSM"; Row = 99;  Cells = @{ "H" = 1058.1;     "I" = 1045.3846;  "J" = 1081.7142;   "K" = 1045.3846;  "L" = 1081.7142;  "M" = 452.6153999999999; "N" = -4077.7142 } },
    @{ Sheet = "BSM"; Row = 107; Cells = @{ "H" = 903.3333;   "I" = 903.3333;   "K" = 903.3333;     "M" = 1016.6667 } },
    @{ Sheet = "BSM"; Row = 134; Cells = @{ "H" = 3486.8;     "I" = 3657.4375;  "K" = 10972.3125;   "M" = -8437.3125 } },

    @{ Sheet = "CRP"; Row = 31;  Cells = @{ "H" = 3063.0466;  "I" = 2489.9375;  "J" = 3402.6667;   "K" = 2489.9375;  "L" = 3402.6667;  "M" = -2194.9375;       "N" = -3992.6667 } },
    @{ Sheet = "CRP"; Row = 34;  Cells = @{ "H" = 3063.0466;  "I" = 2489.9375;  "J" = 3402.6667;   "K" = 2489.9375;  "L" = 3402.6667;  "M" = -2287.9375;       "N" = -3806.6667 } },
    @{ Sheet = "CRP"; Row = 58;  Cells = @{ "H" = 21014.77;   "I" = 1860.6;     "K" = 1860.6;       "M" = -1657.6 } },
    @{ Sheet = "CRP"; Row = 99;  Cells = @{ "H" = 22730704;   "I" = 2886.4;     "K" = 2886.4;       "M" = -1388.4 } },
    @{ Sheet = "CRP"; Row = 126; Cells = @{ "H" = 22730704;   "I" = 2886.4;     "K" = 8659.200000000001; "M" = -6189.200000000001 } },
    @{ Sheet = "CRP"; Row = 134; Cells = @{ "H" = 1030.0555;  "I" = 649.1818;   "K" = 1947.5454;    "M" = 587.4546 } },
    @{ Sheet = "CRP"; Row = 136; Cells = @{ "H" = 21014.77;   "I" = 1860.6;     "K" = 5581.799999999999; "M" = -3031.799999999999 } },

    @{ Sheet = "CUL"; Row = 5;   Cells = @{ "H" = 1264.0555;  "I" = 799.4167;   "J" = 2193.3333;   "K" = 2398.2501;  "L" = 6579.999899999999; "M" = -2286.2501; "N" = -6803.999899999999 } },
    @{ Sheet = "CUL"; Row = 131; Cells = @{ "H" = 754.02;     "I" = 386.25;     "J" = 786;         "K" = 1158.75;    "L" = 2358;       "M" = 3881.25;          "N" = -12438 } },
    @{ Sheet = "CUL"; Row = 132; Cells = @{ "H" = 1000;       "I" = 1000;       "K" = 9000;         "M" = -6470 } },
    @{ Sheet = "CUL"; Row = 135; Cells = @{ "H" = 1264.0555;  "I" = 799.4167;   "J" = 2193.3333;   "K" = 7194.7503;  "L" = 19739.9997; "M" = -4659.7503;       "N" = -24809.9997 } },

    @{ Sheet = "GSM"; Row = 122; Cells = @{ "H" = 95239970;   "I" = 33335254;  "J" = 250001740;   "K" = 100005762; "L" = 750005220;  "M" = -100003312;       "N" = -750010120 } },
    @{ Sheet = "GSM"; Row = 132; Cells = @{ "H" = 28646.422;  "I" = 2001;       "J" = 58252.445;   "K" = 6003;       "L" = 174757.335; "M" = -3473;            "N" = -179817.335 } },

    @{ Sheet = "LTW"; Row = 22;  Cells = @{ "H" = 4728.4287;  "J" = 5499.75;    "L" = 5499.75;      "N" = -6089.75 } },
    @{ Sheet = "LTW"; Row = 27;  Cells = @{ "H" = 4728.4287;  "J" = 5499.75;    "L" = 5499.75;      "N" = -5713.75 } },
    @{ Sheet = "LTW"; Row = 132; Cells = @{ "H" = 671701.7;   "I" = 1340570.2;  "K" = 4021710.6;    "M" = -4019180.6 } },

    @{ Sheet = "WVR"; Row = 70;  Cells = @{ "H" = 20000;      "J" = 20000;      "L" = 20000;        "N" = -20630 } },
    @{ Sheet = "WVR"; Row = 73;  Cells = @{ "H" = 20000;      "J" = 20000;      "L" = 20000;        "N" = -22184 } },
    @{ Sheet = "WVR"; Row = 132; Cells = @{ "H" = 1598.5;     "I" = 1197;       "K" = 3591;          "M" = -1061 } },
    @{ Sheet = "WVR"; Row = 136; Cells = @{ "H" = 38234010;   "I" = 54331224;  "K" = 162993672;     "M" = -162991122 } }
)

foreach ($chg in $changes) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    foreach ($col in $chg.Cells.Keys) {
        $addr = "$col$($chg.Row)"
        # Assigning $null clears/removes the cell (matches rows where a
        # column's value is dropped entirely rather than just updated).
        $ws.Range($addr).Value = $chg.Cells[$col]
    }
}

Write-Output "Applied $($changes.Count) row updates."
